$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared-string bearing value updates ---
# Order matters: it controls the order new shared strings get appended in
# sharedStrings.xml, which must match the target diff.

# Row 2: C2 switches from a raw number 24 to the text "24"
$ws.Range("C2").Value = "24"

# Rows 34-39 (column E): "Op Cost Total Sum ..." -> "Op Cost Total ..."
$ws.Range("E34").Value = "Op Cost Total Current Day"
$ws.Range("E35").Value = "Op Cost Total Current Month"
$ws.Range("E36").Value = "Op Cost Total Current Year"
$ws.Range("E37").Value = "Op Cost Total Previous Day"
$ws.Range("E38").Value = "Op Cost Total Previous Month"
$ws.Range("E39").Value = "Op Cost Total Previous Year"

# Rows 31-33 (column E): "Op kWh Total Counter Current ..." -> "Op kWh Total Current ..."
$ws.Range("E31").Value = "Op kWh Total Current Day"
$ws.Range("E32").Value = "Op kWh Total Current Month"
$ws.Range("E33").Value = "Op kWh Total Current Year"

# Rows 5 & 6 (column H): "1" -> "2"
$ws.Range("H5").Value = "2"
$ws.Range("H6").Value = "2"

# --- Column width / cols cleanup: columns 9-11 -> 9-10 ---
# Column K (11) is empty; deleting it collapses the bestFit width group
# from min=9,max=11 to min=9,max=10 while leaving the real data untouched.
$ws.Range("K1").EntireColumn.Delete()

# --- Selection change ---
$ws.Range("C4").Select() | Out-Null
